$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New data rows to append at the bottom of the timelog (rows 17 and 18)
# Row 17: date/time entry + "click output and brush output" note
$ws.Cells.Item(17, 1).Value = "2/23, 2 hours"
$ws.Cells.Item(17, 2).Value = "worked on click output and brush output"

# Row 18: same date/time entry repeated + "second map" note
$ws.Cells.Item(18, 1).Value = "2/23, 2 hours"
$ws.Cells.Item(18, 2).Value = "Worked on second map which zooms in, adding clicker for second map"

# Copy the formatting (fonts, borders, alignment, wrap) from the prior data row (16)
# onto the two newly added rows so they match the sheet's existing style.
$ws.Range("A16:B16").Copy() | Out-Null
$ws.Range("A17:B18").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Re-apply the row heights seen in the target workbook
$ws.Rows.Item(17).RowHeight = 27.6
$ws.Rows.Item(18).RowHeight = 41.4

# Update the selection to reflect where the user ended up after entering data
$ws.Range("B19").Select() | Out-Null
